$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (as BGR integers expected by the Excel .Color COM property)
$FILL_GRAY   = 10921638   # A6A6A6
$FILL_GOLD   = 4641530    # FAD246
$FILL_WHITE  = 16711679   # FFFFFE
$FILL_LTGRAY = 14277081   # D9D9D9

$BORDER_LTGRAY = 14540253 # DDDDDD
$BORDER_BLACK  = 65536    # 000001
$BORDER_GOLD   = 1952255  # FFC91D

$FONT_NAME = '"Franklin Gothik Book"'

function Set-CellStyle($range, $fillColor, $bold) {
    $range.Font.Name = $FONT_NAME
    $range.Font.Size = 10
    $range.Font.Bold = $bold
    $range.HorizontalAlignment = -4108
    $range.Interior.Pattern = 1
    $range.Interior.Color = $fillColor
}

function Set-Edge($range, $edge, $weight, $color) {
    $b = $range.Borders.Item($edge)
    $b.LineStyle = 1
    $b.Weight = $weight
    $b.Color = $color
}

function Set-BoxBorder($range, $color) {
    Set-Edge $range 7 2 $color
    Set-Edge $range 10 2 $color
    Set-Edge $range 8 2 $color
    Set-Edge $range 9 2 $color
}

# --- Row 114: function signature header (merged B114:C114) ---
$ws.Range("B114").Value = "SimpleRules DoubleValue ratioScore_TRANS(DoubleValue ratio)"
$ws.Range("C114").Value = ""
$hdr = $ws.Range("B114:C114")
$hdr.Merge()
Set-CellStyle $hdr $FILL_WHITE $false
$hdr.Font.Color = $BORDER_BLACK
Set-Edge $hdr 8 2 $BORDER_BLACK
Set-Edge $hdr 9 2 $BORDER_BLACK

# --- Row 115: parameter / RETURN labels ---
$ws.Range("B115").Value = "ratio"
$ws.Range("C115").Value = "RETURN"

$b115 = $ws.Range("B115")
Set-CellStyle $b115 $FILL_GRAY $false
Set-Edge $b115 9 2 $BORDER_BLACK

$c115 = $ws.Range("C115")
Set-CellStyle $c115 $FILL_GOLD $false
Set-Edge $c115 9 4 $BORDER_GOLD

# --- Rows 116-121: ranges and values (alternating white / light-gray) ---
$ranges = @("<=10", "(10 .. 15]", "(15 .. 20]", "(20 .. 30]", "(30 .. 40]", "(40 .. 150]")
$values = @("'0", "'0.7", "'0.8", "'0.9", "'0.95", "'1")
$boldFlags = @($false, $true, $false, $false, $false, $false)
$boldFlagsC = @($true, $false, $false, $false, $false, $false)

for ($i = 0; $i -lt 6; $i++) {
    $row = 116 + $i
    $bCell = $ws.Range("B$row")
    $cCell = $ws.Range("C$row")
    $bCell.Value = $ranges[$i]
    $cCell.Value = $values[$i]
    Set-CellStyle $bCell $FILL_WHITE $boldFlags[$i]
    Set-BoxBorder $bCell $BORDER_LTGRAY
    Set-CellStyle $cCell $FILL_LTGRAY $boldFlagsC[$i]
    Set-BoxBorder $cCell $BORDER_LTGRAY
}

# --- Row 122: last range (closing thick border like the group header) ---
$ws.Range("B122").Value = ">150"
$ws.Range("C122").Value = "'0.95"

$b122 = $ws.Range("B122")
Set-CellStyle $b122 $FILL_WHITE $false
Set-BoxBorder $b122 $BORDER_LTGRAY
Set-Edge $b122 9 2 $BORDER_BLACK

$c122 = $ws.Range("C122")
Set-CellStyle $c122 $FILL_LTGRAY $false
Set-BoxBorder $c122 $BORDER_LTGRAY
Set-Edge $c122 9 4 $BORDER_GOLD

Write-Host "applied ratioScore_TRANS block"
